# Apply cryptos list update (Sun Jun 11 19:49:20 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.068.92"
Set-TextValue $ws.Range("E2") "  +1.49%  "
Set-TextValue $ws.Range("D3") "1.764.56"
Set-TextValue $ws.Range("E3") "  +1.22%  "
Set-TextValue $ws.Range("D4") "0.9995"
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "237.70"
Set-TextValue $ws.Range("D6") "0.9990"
Set-TextValue $ws.Range("E7") "  +4.34%  "
Set-TextValue $ws.Range("D8") "0.2749"
Set-TextValue $ws.Range("E8") "  +3.74%  "
Set-TextValue $ws.Range("D9") "0.06215"
Set-TextValue $ws.Range("E9") "  +1.54%  "
Set-TextValue $ws.Range("D10") "1.771.67"
Set-TextValue $ws.Range("E10") "  +1.63%  "
Set-TextValue $ws.Range("D11") "16.06"
Set-TextValue $ws.Range("E11") "  +5.51%  "
Set-TextValue $ws.Range("D12") "0.07040"
Set-TextValue $ws.Range("E12") "  +1.25%  "
Set-TextValue $ws.Range("D13") "0.6552"
Set-TextValue $ws.Range("E13") "  +9.83%  "
Set-TextValue $ws.Range("D14") "4.519"
Set-TextValue $ws.Range("E14") "  +0.60%  "
Set-TextValue $ws.Range("E15") "  +2.63%  "
Set-TextValue $ws.Range("D16") "0.9987"
Set-TextValue $ws.Range("E16") "  -0.21%  "
Set-TextValue $ws.Range("D17") "0.9991"
Set-TextValue $ws.Range("E17") "  -0.14%  "
Set-TextValue $ws.Range("D18") "26.074.96"
Set-TextValue $ws.Range("E18") "  +1.46%  "
Set-TextValue $ws.Range("D19") "11.74"
Set-TextValue $ws.Range("E19") "  +1.12%  "
Set-TextValue $ws.Range("D20") "0.000006752"
Set-TextValue $ws.Range("E20") "  -0.28%  "
Set-TextValue $ws.Range("D21") "1.996.17"
Set-TextValue $ws.Range("E21") "  +1.42%  "
Set-TextValue $ws.Range("D22") "4.101"
Set-TextValue $ws.Range("E22") "  +1.55%  "
Set-TextValue $ws.Range("D23") "8.443"
Set-TextValue $ws.Range("E23") "  +3.55%  "
Set-TextValue $ws.Range("D24") "5.210"
Set-TextValue $ws.Range("D25") "137.87"
Set-TextValue $ws.Range("E25") "  +0.26%  "
Set-TextValue $ws.Range("D26") "1.485"
Set-TextValue $ws.Range("E26") "  -2.19%  "
Set-TextValue $ws.Range("B27") "LidoDAOToken"
Set-TextValue $ws.Range("C27") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "1.840"
Set-TextValue $ws.Range("E27") "  +0.97%  "
Set-TextValue $ws.Range("B28") "EthereumClassic"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D28") "15.22"
Set-TextValue $ws.Range("E28") "  +1.78%  "
Set-TextValue $ws.Range("D29") "103.09"
Set-TextValue $ws.Range("E29") "  -0.09%  "
Set-TextValue $ws.Range("D30") "0.08429"
Set-TextValue $ws.Range("E30") "  +4.04%  "
Set-TextValue $ws.Range("D31") "3.717"
Set-TextValue $ws.Range("E31") "  -0.97%  "
Set-TextValue $ws.Range("D32") "3.467"
Set-TextValue $ws.Range("E32") "  +0.44%  "
Set-TextValue $ws.Range("D33") "0.04442"
Set-TextValue $ws.Range("E33") "  -1.16%  "
Set-TextValue $ws.Range("D34") "2.651"
Set-TextValue $ws.Range("E34") "  -0.06%  "
Set-TextValue $ws.Range("D35") "1.006"
Set-TextValue $ws.Range("E35") "  +2.50%  "
Set-TextValue $ws.Range("D36") "0.6131"
Set-TextValue $ws.Range("E36") "  +0.56%  "
Set-TextValue $ws.Range("D37") "2.760"
Set-TextValue $ws.Range("E37") "  +4.00%  "
Set-TextValue $ws.Range("D38") "0.01586"
Set-TextValue $ws.Range("E38") "  +2.41%  "
Set-TextValue $ws.Range("D39") "1.985"
Set-TextValue $ws.Range("E39") "  +3.80%  "
Set-TextValue $ws.Range("D40") "1.000"
Set-TextValue $ws.Range("E40") "  +0.06%  "
Set-TextValue $ws.Range("D41") "103.15"
Set-TextValue $ws.Range("E41") "  -0.33%  "
Set-TextValue $ws.Range("D42") "0.3917"
Set-TextValue $ws.Range("E42") "  +3.24%  "
Set-TextValue $ws.Range("D43") "0.7565"
Set-TextValue $ws.Range("E43") "  +4.14%  "
Set-TextValue $ws.Range("D44") "4.980"
Set-TextValue $ws.Range("E44") "  -2.67%  "
Set-TextValue $ws.Range("D45") "0.05505"
Set-TextValue $ws.Range("E45") "  +3.23%  "
Set-TextValue $ws.Range("D46") "6.379"
Set-TextValue $ws.Range("E46") "  +8.46%  "
Set-TextValue $ws.Range("D47") "0.1129"
Set-TextValue $ws.Range("E47") "  +1.52%  "
Set-TextValue $ws.Range("D48") "30.25"
Set-TextValue $ws.Range("E48") "  +0.55%  "
Set-TextValue $ws.Range("D49") "53.07"
Set-TextValue $ws.Range("E49") "  +1.33%  "
Set-TextValue $ws.Range("D50") "0.3472"
Set-TextValue $ws.Range("E50") "  +0.93%  "
Set-TextValue $ws.Range("D51") "1.001"
Set-TextValue $ws.Range("E51") "  +0.30%  "
